$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRANSIENT")

# Insert two new rows before current row 10 (STPMAX row), pushing existing
# rows 10 (STPMAX) and 11 (TAUREF) down to 12 and 13.
$ws.Rows.Item(10).Resize(2).Insert()

# The inserted rows pick up a default style in every column (incl. E); the
# new rows only use columns A-D, so drop the stray E cells entirely.
$ws.Range("E10:E11").Clear()

# Set new shared-string values in the same order they first appear in the
# authoritative edit: MLT_UPPER, MLT_LOWER, "float or str", upper-note, lower-note.
$ws.Range("A10").Value = "MLT_UPPER"
$ws.Range("A11").Value = "MLT_LOWER"
$ws.Range("C10").Value = "float or str"
$ws.Range("C11").Value = "float or str"
$ws.Range("D10").Value = "Multiplier used to tune (increase) the adaptive time step, used if flag IADAPTIME = 1 or IADAPTIME = 2; default to 1.2. Any positive real number or none. If none default value is used. "
$ws.Range("D11").Value = "Multiplier used to tune (decrease) the adaptive time step, used if flag IADAPTIME = 1 or IADAPTIME = 2; default to 0.5. Any positive real number or none. If none default value is used. "

$ws.Range("B10").Value = "-"
$ws.Range("B11").Value = "-"

$ws.Range("A10:A11").VerticalAlignment = -4108  # xlCenter
$ws.Range("B10:C11").HorizontalAlignment = -4108
$ws.Range("B10:C11").VerticalAlignment = -4108
$ws.Range("D10:D11").WrapText = $true
$ws.Rows.Item(10).RowHeight = 29
$ws.Rows.Item(11).RowHeight = 29

# A handful of pre-existing Value cells in column B/C/E switch from the old
# "style 8/9" look to the plain centered "style 2" look used everywhere else
# in the Value column now that 8/9 are being redefined for the new rows.
$ws.Range("E3").HorizontalAlignment = -4108
$ws.Range("E4").HorizontalAlignment = -4108
$ws.Range("B5").HorizontalAlignment = -4108
$ws.Range("C5").HorizontalAlignment = -4108
$ws.Range("E5").HorizontalAlignment = -4108
$ws.Range("E6").HorizontalAlignment = -4108
$ws.Range("E7").HorizontalAlignment = -4108
$ws.Range("E8").HorizontalAlignment = -4108
$ws.Range("E12").HorizontalAlignment = -4108
$ws.Range("E13").HorizontalAlignment = -4108

# Restore the frozen-pane/selection bookkeeping Excel persists for this sheet.
$sheetView = $ws.Application.ActiveWindow
$ws.Range("D16").Select()

Write-Output ("Dimension after insert: " + $ws.UsedRange.Address())
Write-Output ("A10 value: " + $ws.Range("A10").Value)
